# ifo GDP component (GVA) matched-error preprocessing refresh.
#
# A newer ifoCAST run re-matched the quarter-ahead error series against
# the latest vintage, so every existing error in the triangle is
# refreshed, and each row grows by one more observation at its
# "frontier" column (the staircase edge of the triangle).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Sheet, [string]$StartCol, [int]$Row, [string]$EndCol, [double[]]$Values)

    $n = $Values.Count
    $arr2d = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr2d[0, $i] = $Values[$i]
    }
    $Sheet.Range("$StartCol$Row" + ":" + "$EndCol$Row").Value = $arr2d
}

Set-RowValues $ws "B" 2  "I" @(-0.1140116980535687, 1.134027854149689, 0.5293419069744123, 1.034669590998859, 1.285355063281015, 1.444240868935941, 0.3125477243455349, -0.5548516401491572)
Set-RowValues $ws "B" 3  "H" @(0.52162968558209, 0.6031613697014773, 1.324783753290712, 1.4345885005637, 1.547161051919687, 0.4466123193550214, -0.4809009009009035)
Set-RowValues $ws "B" 4  "K" @(-0.5273222390317647, -0.2552162467092876, -0.005411499436300538, 0.7074136955291455, 0.3090111142775654, -0.4209009009009035, 1.271278195488724, 0.6220589070934801, 0.8193626707131937, 0.696085300837782)
Set-RowValues $ws "B" 5  "J" @(-0.3652162467092876, 0.004588500563699471, 0.7769453239451327, 0.3807516038255184, -0.3988174318051805, 1.180575399060113, 0.5077756378413858, 0.6787354624760854, 0.52058802116979)
Set-RowValues $ws "B" 6  "I" @(0.1867105431484077, 0.7722779613039183, 0.2843558036209001, -0.3395270793998535, 1.287674180139613, 0.5677696640183092, 0.6675525209407469, 0.5067171635883239)
Set-RowValues $ws "B" 7  "H" @(0.8280615287077258, 0.009011114277565485, -0.4809009009009035, 1.261278195488724, 0.6020589070934801, 0.7093626707131937, 0.606085300837782)
Set-RowValues $ws "B" 8  "J" @(-0.3509888857224345, -0.5809009009009034, 1.241278195488724, 0.7032066015985186, 0.8825435650546625, 0.6567688549060375, -0.1837499599687909, 0.3538112635634896, 0.1485032540344368)
Set-RowValues $ws "B" 9  "I" @(-0.8569009009009034, 1.055249540297363, 0.5422953553339147, 0.8185372705227552, 0.7162146587346183, -0.115966644056028, 0.3949757811758082, 0.1734537503564907)
Set-RowValues $ws "B" 10 "H" @(0.6787681550252201, 0.3564167683417685, 0.8907192198060601, 0.8080775651369163, -0.0717095860245372, 0.4455085207778247, 0.2422520263583712)
Set-RowValues $ws "B" 11 "G" @(0.4003114290510102, 0.8093626707131937, 0.796085300837782, -0.05551312649165269, 0.4849757811758083, 0.2534537503564908)
Set-RowValues $ws "B" 12 "F" @(0.33218714548825, 0.486085300837782, -0.2025131264916527, 0.5149757811758082, 0.381103329907261)
Set-RowValues $ws "B" 13 "E" @(0.4903150813867297, -0.3157490177063643, 0.2454495103302783, 0.04235042473292953)
Set-RowValues $ws "B" 14 "D" @(-0.2760759690770392, 0.2958161606567877, 0.07961008106920435)
Set-RowValues $ws "B" 15 "C" @(0.1421546153588515, 0.02893023050567838)
Set-RowValues $ws "B" 16 "B" @(0.02940328597706714)
